$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 161063
$ws.Range("C4").Value = 152084
$ws.Range("C5").Value = 8979
$ws.Range("C7").Value = 5.57
$ws.Range("C8").Value = 64.5
